# Insert two new data rows into the weekly price table.
# The new rows are inserted right before the current row 338, shifting the
# existing rows 338..446 down to 340..448 (the two rows that fall off the
# bottom become the new rows 447 and 448).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 338:339 — everything from the old row 338 onward
# shifts down by two rows (formatting/number format of column D is carried
# along automatically by Excel's row insert).
$ws.Rows("338:339").Insert()

# --- New row 338 ---------------------------------------------------------
$ws.Range("A338").Value = 4
$ws.Range("B338").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C338").Value = "Los Lagos"
$ws.Range("D338").Value = 45215
$ws.Range("E338").Value = 10
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100101
$ws.Range("H338").Value = "Berries"
$ws.Range("I338").Value = 100112025
$ws.Range("J338").Value = "Frutilla"
$ws.Range("K338").Value = "Sin especificar"
$ws.Range("L338").Value = "Especial"
$ws.Range("M338").Value = 300
$ws.Range("N338").Value = 13000
$ws.Range("O338").Value = 13000
$ws.Range("P338").Value = 13000
$ws.Range("Q338").Value = "`$/bandeja 7 kilos"
$ws.Range("R338").Value = "Provincia de Melipilla"
$ws.Range("S338").Value = 1857
$ws.Range("T338").Value = 7

# --- New row 339 ---------------------------------------------------------
$ws.Range("A339").Value = 4
$ws.Range("B339").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C339").Value = "Los Lagos"
$ws.Range("D339").Value = 45215
$ws.Range("E339").Value = 10
$ws.Range("F339").Value = "Fruta"
$ws.Range("G339").Value = 100101
$ws.Range("H339").Value = "Berries"
$ws.Range("I339").Value = 100112025
$ws.Range("J339").Value = "Frutilla"
$ws.Range("K339").Value = "Sin especificar"
$ws.Range("L339").Value = "Primera"
$ws.Range("M339").Value = 300
$ws.Range("N339").Value = 10000
$ws.Range("O339").Value = 10000
$ws.Range("P339").Value = 10000
$ws.Range("Q339").Value = "`$/bandeja 7 kilos"
$ws.Range("R339").Value = "Provincia de Melipilla"
$ws.Range("S339").Value = 1429
$ws.Range("T339").Value = 7
